# Error Calculations and Plots
# Applies the "missing_data" edits: two rows are removed entirely (RM 232 and
# SC 92), and a number of individual A-F cells are blanked out or filled back
# in across the remaining rows, simulating a different random selection of
# missing values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two deleted records -----------------------------------
# Delete the higher row index first so the lower one's position is not
# affected by the shift.
$ws.Rows.Item(28).Delete()   # "SC 92"
$ws.Rows.Item(26).Delete()   # "RM 232"

# --- Blank out cells that became missing -------------------------------
$ws.Range("E2").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("F17").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("C20").ClearContents()
$ws.Range("F25").ClearContents()
$ws.Range("F26").ClearContents()
$ws.Range("C27").ClearContents()
$ws.Range("F27").ClearContents()
$ws.Range("B28").ClearContents()
$ws.Range("E28").ClearContents()
$ws.Range("F28").ClearContents()
$ws.Range("B29").ClearContents()
$ws.Range("B32").ClearContents()

# --- Fill in cells that now have a known value --------------------------
$ws.Range("E5").Value = -5
$ws.Range("C6").Value = 15.1
$ws.Range("E6").Value = -5.7
$ws.Range("F6").Value = 16.43
$ws.Range("F11").Value = 17.65
$ws.Range("C12").Value = 12.5
$ws.Range("F14").Value = 17.76
$ws.Range("C17").Value = 11.2
$ws.Range("C18").Value = 11.5
$ws.Range("F19").Value = 17.81
$ws.Range("F21").Value = 16.58
$ws.Range("F22").Value = 16.81
$ws.Range("C23").Value = 12.2
$ws.Range("E24").Value = -8.1
$ws.Range("B27").Value = -20.4
$ws.Range("B30").Value = -19.7
$ws.Range("E30").Value = -5.7
$ws.Range("F31").Value = 17.18
